# Update EncounterTypes sheet: reshuffle row 3, replace row 4, append new
# "Strong heat / thirst" encounter rows, and highlight the affected cells
# with a yellow fill.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Arriving at Port / fraud / merchants selling supplies ---------
$ws.Range("B3").Value = "Arriving at Port of Krasnovdsk"
$ws.Range("C3").Value = "Turcoman merchant committing fraud"
$ws.Range("D3").Value = "Turcoman merchants selling supplies"

# --- Row 4 onward: drop "Soldiers dying due to diseases", shift the rest up
$ws.Range("C4").Value = "Oasis mirage"
$ws.Range("C5").Value = "Strong heat causing thirst - mild"
$ws.Range("C6").Value = "Strong heat causing thirst - dangerous"
$ws.Range("C7").Value = "Strong heat causing thirst - extremely dangerous"
$ws.Range("C8").Value = "Low supply, Major Frankenburg in verge of death"
$ws.Range("C9").Value = "Turcoman dying of thirst"

# --- Column width for column D ---------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 35.5

# --- Highlight cells with a yellow fill (all data cells except C3 and C7) -
# Excel's Interior.Color takes a BGR-packed long; 255,255,0 = yellow.
$yellow = 255 + (255 * 256) + (0 * 65536)

$highlightRanges = @("B2", "C2", "D2", "B3", "D3", "C4", "C5", "C6", "C8", "C9")
foreach ($addr in $highlightRanges) {
    $ws.Range($addr).Interior.Color = $yellow
}

# --- Update selection to match the target state -----------------------------
$ws.Range("D2").Select()
